$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet SCD0182 -> SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-197" to "SCD0011-013"
$ws.Range("B2").Value = "SCD0011-013"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.67

# Move the active selection from O2 to B3 (and drop the E1 scroll anchor)
$ws.Range("B3").Select()
